$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(10, 8).Value = 3252  # H10
$ws.Cells.Item(10, 10).Value = 5000  # J10
$ws.Cells.Item(10, 12).Value = 5000  # L10
$ws.Cells.Item(10, 14).Value = -5586  # N10
$ws.Cells.Item(12, 8).Value = 0  # H12
$ws.Cells.Item(12, 9).Value = 0  # I12
$ws.Cells.Item(12, 10).Value = 0  # J12
$ws.Cells.Item(12, 11).Value = 0  # K12
$ws.Cells.Item(12, 12).Value = 0  # L12
$ws.Cells.Item(12, 13).ClearContents()  # M12
$ws.Cells.Item(12, 14).ClearContents()  # N12
$ws.Cells.Item(19, 8).Value = 1953.8334  # H19
$ws.Cells.Item(19, 9).Value = 779.5714  # I19
$ws.Cells.Item(19, 10).Value = 3597.8  # J19
$ws.Cells.Item(19, 11).Value = 779.5714  # K19
$ws.Cells.Item(19, 12).Value = 3597.8  # L19
$ws.Cells.Item(19, 13).Value = -604.5714  # M19
$ws.Cells.Item(19, 14).Value = -3947.8  # N19
$ws.Cells.Item(32, 9).Value = 900  # I32
$ws.Cells.Item(32, 10).Value = 1262.25  # J32
$ws.Cells.Item(32, 11).Value = 900  # K32
$ws.Cells.Item(32, 12).Value = 1262.25  # L32
$ws.Cells.Item(32, 13).Value = -574  # M32
$ws.Cells.Item(32, 14).Value = -1914.25  # N32
$ws.Cells.Item(38, 8).Value = 762.375  # H38
$ws.Cells.Item(38, 9).Value = 762.375  # I38
$ws.Cells.Item(38, 11).Value = 2287.125  # K38
$ws.Cells.Item(38, 13).Value = -1915.125  # M38
$ws.Cells.Item(42, 8).Value = 1458.375  # H42
$ws.Cells.Item(42, 9).Value = 2569.5  # I42
$ws.Cells.Item(42, 10).Value = 347.25  # J42
$ws.Cells.Item(42, 11).Value = 7708.5  # K42
$ws.Cells.Item(42, 12).Value = 1041.75  # L42
$ws.Cells.Item(42, 13).Value = -7478.5  # M42
$ws.Cells.Item(42, 14).Value = -1501.75  # N42
$ws.Cells.Item(76, 8).Value = 1806122.8  # H76
$ws.Cells.Item(76, 9).Value = 2930752.5  # I76
$ws.Cells.Item(76, 11).Value = 2930752.5  # K76
$ws.Cells.Item(76, 13).Value = -2930437.5  # M76
$ws.Cells.Item(79, 8).Value = 1806122.8  # H79
$ws.Cells.Item(79, 9).Value = 2930752.5  # I79
$ws.Cells.Item(79, 11).Value = 2930752.5  # K79
$ws.Cells.Item(79, 13).Value = -2929660.5  # M79
$ws.Cells.Item(132, 8).Value = 1120.5807  # H132
$ws.Cells.Item(132, 9).Value = 1128.931  # I132
$ws.Cells.Item(132, 10).Value = 999.5  # J132
$ws.Cells.Item(132, 11).Value = 3386.793  # K132
$ws.Cells.Item(132, 12).Value = 2998.5  # L132
$ws.Cells.Item(132, 13).Value = -856.7930000000001  # M132
$ws.Cells.Item(132, 14).Value = -8058.5  # N132
$ws.Cells.Item(137, 8).Value = 1889.091  # H137
$ws.Cells.Item(137, 9).Value = 1231.8889  # I137
$ws.Cells.Item(137, 11).Value = 3695.6667  # K137
$ws.Cells.Item(137, 13).Value = -1145.6667  # M137
$ws.Cells.Item(140, 8).Value = 52259.617  # H140
$ws.Cells.Item(140, 10).Value = 52259.617  # J140
$ws.Cells.Item(140, 12).Value = 52259.617  # L140
$ws.Cells.Item(140, 14).Value = -62619.617  # N140

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4889.3335  # H32
$ws.Cells.Item(32, 9).Value = 3305.6206  # I32
$ws.Cells.Item(32, 10).Value = 7759.8125  # J32
$ws.Cells.Item(32, 11).Value = 3305.6206  # K32
$ws.Cells.Item(32, 12).Value = 7759.8125  # L32
$ws.Cells.Item(32, 13).Value = -3018.6206  # M32
$ws.Cells.Item(32, 14).Value = -8333.8125  # N32
$ws.Cells.Item(63, 8).Value = 3499.4  # H63
$ws.Cells.Item(63, 9).Value = 3499.4  # I63
$ws.Cells.Item(63, 11).Value = 3499.4  # K63
$ws.Cells.Item(63, 13).Value = -2813.4  # M63
$ws.Cells.Item(66, 8).Value = 3499.4  # H66
$ws.Cells.Item(66, 9).Value = 3499.4  # I66
$ws.Cells.Item(66, 11).Value = 17497  # K66
$ws.Cells.Item(66, 13).Value = -14065  # M66
$ws.Cells.Item(74, 8).Value = 866.6667  # H74
$ws.Cells.Item(74, 10).Value = 1406.5  # J74
$ws.Cells.Item(74, 12).Value = 1406.5  # L74
$ws.Cells.Item(74, 14).Value = -3154.5  # N74
$ws.Cells.Item(77, 8).Value = 866.6667  # H77
$ws.Cells.Item(77, 10).Value = 1406.5  # J77
$ws.Cells.Item(77, 12).Value = 7032.5  # L77
$ws.Cells.Item(77, 14).Value = -15768.5  # N77
$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 10).Value = 0  # J133
$ws.Cells.Item(133, 12).Value = 0  # L133
$ws.Cells.Item(133, 14).ClearContents()  # N133

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2151.75  # H20
$ws.Cells.Item(20, 9).Value = 1922  # I20
$ws.Cells.Item(20, 10).Value = 3300.5  # J20
$ws.Cells.Item(20, 11).Value = 1922  # K20
$ws.Cells.Item(20, 12).Value = 3300.5  # L20
$ws.Cells.Item(20, 13).Value = -1675  # M20
$ws.Cells.Item(20, 14).Value = -3794.5  # N20
$ws.Cells.Item(80, 8).Value = 8595.166999999999  # H80
$ws.Cells.Item(80, 10).Value = 10304.9  # J80
$ws.Cells.Item(80, 12).Value = 10304.9  # L80
$ws.Cells.Item(80, 14).Value = -12300.9  # N80
$ws.Cells.Item(83, 8).Value = 8595.166999999999  # H83
$ws.Cells.Item(83, 10).Value = 10304.9  # J83
$ws.Cells.Item(83, 12).Value = 51524.5  # L83
$ws.Cells.Item(83, 14).Value = -61508.5  # N83
$ws.Cells.Item(94, 8).Value = 748.06665  # H94
$ws.Cells.Item(94, 9).Value = 678.61536  # I94
$ws.Cells.Item(94, 10).Value = 1199.5  # J94
$ws.Cells.Item(94, 11).Value = 678.61536  # K94
$ws.Cells.Item(94, 12).Value = 1199.5  # L94
$ws.Cells.Item(94, 13).Value = -227.61536  # M94
$ws.Cells.Item(94, 14).Value = -2101.5  # N94
$ws.Cells.Item(105, 8).Value = 1902  # H105
$ws.Cells.Item(105, 9).Value = 2050  # I105
$ws.Cells.Item(105, 11).Value = 2050  # K105
$ws.Cells.Item(105, 13).Value = -303  # M105
$ws.Cells.Item(119, 8).Value = 40000  # H119
$ws.Cells.Item(119, 10).Value = 40000  # J119
$ws.Cells.Item(119, 12).Value = 40000  # L119
$ws.Cells.Item(119, 14).Value = -49676  # N119

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 7500  # H29
$ws.Cells.Item(29, 10).Value = 7500  # J29
$ws.Cells.Item(29, 12).Value = 7500  # L29
$ws.Cells.Item(29, 14).Value = -8086  # N29
$ws.Cells.Item(132, 8).Value = 2196.2222  # H132
$ws.Cells.Item(132, 9).Value = 1588.4445  # I132
$ws.Cells.Item(132, 10).Value = 2804  # J132
$ws.Cells.Item(132, 11).Value = 4765.333500000001  # K132
$ws.Cells.Item(132, 12).Value = 8412  # L132
$ws.Cells.Item(132, 13).Value = -2235.333500000001  # M132
$ws.Cells.Item(132, 14).Value = -13472  # N132
$ws.Cells.Item(134, 8).Value = 867.8  # H134
$ws.Cells.Item(134, 9).Value = 885.1667  # I134
$ws.Cells.Item(134, 11).Value = 2655.5001  # K134
$ws.Cells.Item(134, 13).Value = -120.5001000000002  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 734.3333  # H5
$ws.Cells.Item(5, 10).Value = 1002.5  # J5
$ws.Cells.Item(5, 12).Value = 3007.5  # L5
$ws.Cells.Item(5, 14).Value = -3231.5  # N5
$ws.Cells.Item(51, 8).Value = 0  # H51
$ws.Cells.Item(51, 9).Value = 0  # I51
$ws.Cells.Item(51, 11).Value = 0  # K51
$ws.Cells.Item(51, 13).ClearContents()  # M51
$ws.Cells.Item(56, 8).Value = 15319.5  # H56
$ws.Cells.Item(56, 9).Value = 15319.5  # I56
$ws.Cells.Item(56, 11).Value = 15319.5  # K56
$ws.Cells.Item(56, 13).Value = -14789.5  # M56
$ws.Cells.Item(122, 8).Value = 609.75  # H122
$ws.Cells.Item(122, 10).Value = 694.5  # J122
$ws.Cells.Item(122, 12).Value = 6250.5  # L122
$ws.Cells.Item(122, 14).Value = -11150.5  # N122
$ws.Cells.Item(131, 8).Value = 776.3  # H131
$ws.Cells.Item(131, 10).Value = 812.86957  # J131
$ws.Cells.Item(131, 12).Value = 2438.60871  # L131
$ws.Cells.Item(131, 14).Value = -12518.60871  # N131
$ws.Cells.Item(135, 8).Value = 734.3333  # H135
$ws.Cells.Item(135, 10).Value = 1002.5  # J135
$ws.Cells.Item(135, 12).Value = 9022.5  # L135
$ws.Cells.Item(135, 14).Value = -14092.5  # N135
$ws.Cells.Item(136, 8).Value = 3006.25  # H136
$ws.Cells.Item(136, 9).Value = 3006.25  # I136
$ws.Cells.Item(136, 11).Value = 9018.75  # K136
$ws.Cells.Item(136, 13).Value = -3918.75  # M136
$ws.Cells.Item(137, 8).Value = 1991.4286  # H137
$ws.Cells.Item(137, 9).Value = 2810  # I137
$ws.Cells.Item(137, 10).Value = 900  # J137
$ws.Cells.Item(137, 11).Value = 8430  # K137
$ws.Cells.Item(137, 12).Value = 2700  # L137
$ws.Cells.Item(137, 13).Value = -3330  # M137
$ws.Cells.Item(137, 14).Value = -12900  # N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5131.846  # H70
$ws.Cells.Item(70, 10).Value = 4263.4  # J70
$ws.Cells.Item(70, 12).Value = 4263.4  # L70
$ws.Cells.Item(70, 14).Value = -4803.4  # N70
$ws.Cells.Item(73, 8).Value = 5131.846  # H73
$ws.Cells.Item(73, 10).Value = 4263.4  # J73
$ws.Cells.Item(73, 12).Value = 4263.4  # L73
$ws.Cells.Item(73, 14).Value = -6135.4  # N73
$ws.Cells.Item(104, 8).Value = 50000  # H104
$ws.Cells.Item(104, 10).Value = 50000  # J104
$ws.Cells.Item(104, 12).Value = 50000  # L104
$ws.Cells.Item(104, 14).Value = -56988  # N104
$ws.Cells.Item(113, 8).Value = 1657.6666  # H113
$ws.Cells.Item(113, 9).Value = 1721  # I113
$ws.Cells.Item(113, 11).Value = 1721  # K113
$ws.Cells.Item(113, 13).Value = 449  # M113

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5348.7144  # H132
$ws.Cells.Item(132, 9).Value = 1180.1177  # I132
$ws.Cells.Item(132, 11).Value = 3540.3531  # K132
$ws.Cells.Item(132, 13).Value = -1010.3531  # M132
